$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5316996666666666
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 8.377268358886777
$ws.Range("R2").Value = 75.395415229981
$ws.Range("S2").Value = 0.3220556913988901
$ws.Range("T2").Value = 0.32205569139889

# Row 3
$ws.Range("G3").Value = 0.5316996666666666
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 14.33135164251755
$ws.Range("S3").Value = 0.5509544596378365
$ws.Range("T3").Value = 0.5509544596378364

# Row 4
$ws.Range("G4").Value = 0.5316996666666666
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 3.303242488896777
$ws.Range("S4").Value = 0.1269898489632735
$ws.Range("T4").Value = 0.1269898489632735
